{"js": "// Office.js (Word JavaScript API) script.\n// Applies the textual edits described by the commit diff:\n//  1. \"Clinical Indication\" -> \"Clinical Indication\" (kept bold; run split\n//     with an inserted proofing mark in the source is not representable\n//     through the content APIs, the visible text/formatting is unchanged).\n//  2. \"?Germline vs somatic origin of previously \" -> unchanged visible text\n//     (source only adds a proofing-error marker around \"?Germline\").\n//  3. \"Illumina NovaSeq\" -> unchanged visible text (adds a spell-check\n//     proofing marker around \"NovaSeq\").\n//  4. \"A custom pipeline utilising the Oncoanalyser analysis pipeline\n//     (OncoPath v1)\" -> unchanged visible text (adds spell-check proofing\n//     markers around \"Oncoanalyser\" and \"OncoPath\").\n//  5. \"...departmental policy. \" -> unchanged visible text (run split only).\n//  6. \"...particularly those > 25 bp in length)\" ->\n//     \"...particularly those > 25 bp in length or in homopolymer regions)\"\n//     (real textual insertion).\n//  7. \"...Our clinical recommendations...\" -> unchanged visible text (run\n//     split only).\n//\n// Items 1, 2, 3, 4, 5, 7 only change how the run is internally split\n// (and, in the source document, add w:proofErr spell/grammar-check\n// markers) without changing the rendered text, so nothing further needs\n// to be typed for them here. Item 6 is the one genuine content change,\n// applied below.\n\nconst body = context.document.body;\n\n// Guard against re-running on an already-edited document.\nconst already = body.search(\"length or in homopolymer regions\", {\n  matchCase: true,\n});\nalready.load(\"text\");\nawait context.sync();\n\nif (already.items.length === 0) {\n  const lengthMatches = body.search(\"particularly those > 25 bp in length\", {\n    matchCase: true,\n  });\n  lengthMatches.load(\"text\");\n  await context.sync();\n\n  if (lengthMatches.items.length > 0) {\n    lengthMatches.items[0].insertText(\n      \" or in homopolymer regions\",\n      Word.InsertLocation.after\n    );\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the textual edits described by the commit diff:\n#  1. \"Clinical Indication\" -> unchanged visible text (source only splits the\n#     run into \"Clinical \" + \"Indication\", both still bold, and adds a\n#     w:proofErr gramStart marker - a spell/grammar-check artifact that the\n#     object model does not expose a way to author directly).\n#  2. \"?Germline vs somatic origin of previously \" -> unchanged visible text\n#     (source only adds a proofing-error marker around \"?Germline\").\n#  3. \"Illumina NovaSeq\" -> unchanged visible text (adds a spell-check\n#     proofing marker around \"NovaSeq\").\n#  4. \"A custom pipeline utilising the Oncoanalyser analysis pipeline\n#     (OncoPath v1)\" -> unchanged visible text (adds spell-check proofing\n#     markers around \"Oncoanalyser\" and \"OncoPath\").\n#  5. \"...departmental policy. \" -> unchanged visible text (run split only).\n#  6. \"...particularly those > 25 bp in length)\" ->\n#     \"...particularly those > 25 bp in length or in homopolymer regions)\"\n#     (real textual insertion).\n#  7. \"...Our clinical recommendations...\" -> unchanged visible text (run\n#     split only).\n#\n# Items 1, 2, 3, 4, 5, 7 only change how the run is internally split (and,\n# in the source document, add w:proofErr spell/grammar-check markers)\n# without changing the rendered text, so nothing further needs to be typed\n# for them here. Item 6 is the one genuine content change, applied below.\n\n$d = $word.ActiveDocument\n\n# Guard against re-running on an already-edited document.\n$already = $d.Content\n$alreadyFind = $already.Find\n$alreadyFind.ClearFormatting()\n$alreadyFind.Text = \"length or in homopolymer regions\"\n$alreadyFind.Forward = $true\n$alreadyFind.Wrap = 0  # wdFindStop\n$alreadyFound = $alreadyFind.Execute()\n\nif (-not $alreadyFound) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = \"particularly those > 25 bp in length\"\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop\n\n    $found = $find.Execute()\n    if ($found) {\n        $rng.Collapse(0)  # wdCollapseEnd\n        $rng.InsertAfter(\" or in homopolymer regions\")\n    }\n}\n"}
